# Auto update: 2025-12-03 10:00:39
# Refresh the daily 국장 조선(shipbuilding) screening table with the latest
# pull: prices/indicators recompute for rows 3-5, and the SamsungHvyInd /
# Hanwha Ocean rows swap places (row 4 <-> row 5) after re-sorting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (HD HYUNDAI MIPO / 010620.KS) - only the MACRO_SCORE refresh changes
$ws.Range("N2").Value = 66.09241856096124

# Row 3 (HDKSOE / 009540.KS)
$ws.Range("D3").Value = 419500
$ws.Range("E3").Value = 42.9
$ws.Range("F3").Value = -1.99
$ws.Range("H3").Value = 66
$ws.Range("I3").Value = 66
$ws.Range("J3").Value = 73
$ws.Range("K3").Value = 52.2
$ws.Range("N3").Value = 66.09241856096124

# Row 4 now holds SamsungHvyInd / 010140.KS (was Hanwha Ocean / 042660.KS)
$ws.Range("B4").Value = "SamsungHvyInd"
$ws.Range("C4").Value = "010140.KS"
$ws.Range("D4").Value = 24575
$ws.Range("E4").Value = 34.2
$ws.Range("F4").Value = -3.06
$ws.Range("H4").Value = 66
$ws.Range("I4").Value = 73
$ws.Range("J4").Value = 96
$ws.Range("K4").Value = 52
$ws.Range("N4").Value = 66.09241856096124

# Row 5 now holds Hanwha Ocean / 042660.KS (was SamsungHvyInd / 010140.KS)
$ws.Range("B5").Value = "Hanwha Ocean"
$ws.Range("C5").Value = "042660.KS"
$ws.Range("D5").Value = 107500
$ws.Range("E5").Value = 19.3
$ws.Range("F5").Value = -5.62
$ws.Range("H5").Value = 63
$ws.Range("I5").Value = 70
$ws.Range("J5").Value = 76
$ws.Range("K5").Value = 50.8
$ws.Range("N5").Value = 66.09241856096124
